# The sheet had an empty row 2 between the header (row 1) and the data
# rows (rows 3 and 4). Deleting that empty row shifts the data rows up
# by one (row 3 -> row 2, row 4 -> row 3), matching the target layout
# where the data now starts immediately after the header with no gap.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(2).Delete()
